$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rarres2"
$ws.Range("C2").Value = "Cmklr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.248266
$ws.Range("H2").Value = 3.744798
$ws.Range("I2").Value = 0.008714053836811378
$ws.Range("J2").Value = 0.00871405383681138
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.041769666666667
$ws.Range("N2").Value = 6.125309
$ws.Range("O2").Value = 0.02341906427171577
$ws.Range("P2").Value = 0.02341906427171577
$ws.Range("Q2").Value = 2.548671654731334
$ws.Range("R2").Value = 22.938044892582
$ws.Range("S2").Value = 0.0002040749868714771
$ws.Range("T2").Value = 0.0002040749868714771

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rarres2"
$ws.Range("C3").Value = "Cmklr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.248266
$ws.Range("H3").Value = 3.744798
$ws.Range("I3").Value = 0.008714053836811378
$ws.Range("J3").Value = 0.00871405383681138
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 84.83061466666666
$ws.Range("N3").Value = 254.491844
$ws.Range("O3").Value = 0.9730057457123328
$ws.Range("P3").Value = 0.973005745712333
$ws.Range("Q3").Value = 105.8911720475013
$ws.Range("R3").Value = 953.020548427512
$ws.Range("S3").Value = 0.00847882445166407
$ws.Range("T3").Value = 0.008478824451664074

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Rarres2"
$ws.Range("C4").Value = "Cmklr1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.248266
$ws.Range("H4").Value = 3.744798
$ws.Range("I4").Value = 0.008714053836811378
$ws.Range("J4").Value = 0.00871405383681138
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3116996666666667
$ws.Range("N4").Value = 0.9350989999999999
$ws.Range("O4").Value = 0.003575190015951382
$ws.Range("P4").Value = 0.003575190015951383
$ws.Range("Q4").Value = 0.3890840961113333
$ws.Range("R4").Value = 3.501756865002
$ws.Range("S4").Value = 0.00003115439827583087
$ws.Range("T4").Value = 0.00003115439827583089

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rarres2"
$ws.Range("C5").Value = "Cmklr1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 103.907654
$ws.Range("H5").Value = 311.722962
$ws.Range("I5").Value = 0.7253717484997341
$ws.Range("J5").Value = 0.7253717484997342
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 2.041769666666667
$ws.Range("N5").Value = 6.125309
$ws.Range("O5").Value = 0.02341906427171577
$ws.Range("P5").Value = 0.02341906427171577
$ws.Range("Q5").Value = 212.1554960716953
$ws.Range("R5").Value = 1909.399464645258
$ws.Range("S5").Value = 0.01698752759900212
$ws.Range("T5").Value = 0.01698752759900212

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Rarres2"
$ws.Range("C6").Value = "Cmklr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 103.907654
$ws.Range("H6").Value = 311.722962
$ws.Range("I6").Value = 0.7253717484997341
$ws.Range("J6").Value = 0.7253717484997342
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 84.83061466666666
$ws.Range("N6").Value = 254.491844
$ws.Range("O6").Value = 0.9730057457123328
$ws.Range("P6").Value = 0.973005745712333
$ws.Range("Q6").Value = 8814.550157391324
$ws.Range("R6").Value = 79330.95141652193
$ws.Range("S6").Value = 0.7057908790676425
$ws.Range("T6").Value = 0.7057908790676427

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Rarres2"
$ws.Range("C7").Value = "Cmklr1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 103.907654
$ws.Range("H7").Value = 311.722962
$ws.Range("I7").Value = 0.7253717484997341
$ws.Range("J7").Value = 0.7253717484997342
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.3116996666666667
$ws.Range("N7").Value = 0.9350989999999999
$ws.Range("O7").Value = 0.003575190015951382
$ws.Range("P7").Value = 0.003575190015951383
$ws.Range("Q7").Value = 32.38798111591533
$ws.Range("R7").Value = 291.4918300432379
$ws.Range("S7").Value = 0.002593341833089446
$ws.Range("T7").Value = 0.002593341833089447

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Rarres2"
$ws.Range("C8").Value = "Cmklr1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 38.091531
$ws.Range("H8").Value = 114.274593
$ws.Range("I8").Value = 0.2659141976634544
$ws.Range("J8").Value = 0.2659141976634544
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 2.041769666666667
$ws.Range("N8").Value = 6.125309
$ws.Range("O8").Value = 0.02341906427171577
$ws.Range("P8").Value = 0.02341906427171577
$ws.Range("Q8").Value = 77.774132552693
$ws.Range("R8").Value = 699.967192974237
$ws.Range("S8").Value = 0.00622746168584217
$ws.Range("T8").Value = 0.006227461685842171

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Rarres2"
$ws.Range("C9").Value = "Cmklr1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 38.091531
$ws.Range("H9").Value = 114.274593
$ws.Range("I9").Value = 0.2659141976634544
$ws.Range("J9").Value = 0.2659141976634544
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 84.83061466666666
$ws.Range("N9").Value = 254.491844
$ws.Range("O9").Value = 0.9730057457123328
$ws.Range("P9").Value = 0.973005745712333
$ws.Range("Q9").Value = 3231.327988324388
$ws.Range("R9").Value = 29081.95189491949
$ws.Range("S9").Value = 0.2587360421930262
$ws.Range("T9").Value = 0.2587360421930262

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Rarres2"
$ws.Range("C10").Value = "Cmklr1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 38.091531
$ws.Range("H10").Value = 114.274593
$ws.Range("I10").Value = 0.2659141976634544
$ws.Range("J10").Value = 0.2659141976634544
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.3116996666666667
$ws.Range("N10").Value = 0.9350989999999999
$ws.Range("O10").Value = 0.003575190015951382
$ws.Range("P10").Value = 0.003575190015951383
$ws.Range("Q10").Value = 11.873117515523
$ws.Range("R10").Value = 106.858057639707
$ws.Range("S10").Value = 0.0009506937845861047
$ws.Range("T10").Value = 0.0009506937845861048

